$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card20")

# N1 header text loses its trailing space: "Correction " -> "Correction"
$ws.Range("N1").Value = "Correction"

# New column O: header "Serviced by " (with trailing space), copying N1's
# header formatting (bold, centered/top aligned, thin border) so it matches
# the rest of the header row instead of getting a brand-new style.
$ws.Range("O1").Value = "Serviced by "
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)

# Data rows 2-12: N gets "nan" (matching the rest of the row) and the new
# O column is created but left blank. Touching .Style on the new O cells
# (re-applying the default "Normal" style) is enough to make Excel keep
# the now-existing-but-empty cell in the sheet instead of dropping it.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 14).Value = "nan"
    $ws.Cells.Item($r, 15).Style = "Normal"
}
